$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coke")
$ws.Range("Z20").Value = "testA"
$ws.Range("Z20").Style = "Normal 2"
